# Update the "dSF" (column F) values for the pitch-data rows.
# Source data was repulled; only column F changed for rows 2,3,4,5,7 (row 6 unchanged).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 4
$ws.Range("F3").Value = 6
$ws.Range("F4").Value = 10
$ws.Range("F5").Value = 4
$ws.Range("F7").Value = 3
